$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.625.13"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "1.753.00"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'324.24"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4601"
$ws.Range("E7").Value = "  +8.41%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "'0.07522"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").Value = "'42.25"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").Value = "'1.098"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "'20.78"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").Value = "'6.018"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "'7.111"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "1.753.49"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "'92.43"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "'0.00001068"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "'0.06420"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'16.77"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").Value = "'5.826"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "27.656.08"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'2.116"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'164.16"
$ws.Range("E26").Value = "  +4.15%  "
$ws.Range("D27").Value = "'20.44"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("D28").Value = "1.954.77"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "'2.094"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'126.97"
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("D31").Value = "'1.075"
$ws.Range("E31").Value = "  -6.92%  "
$ws.Range("D32").Value = "'0.09219"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("D33").Value = "'3.672"
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").Value = "'5.532"
$ws.Range("E34").Value = "  -1.61%  "
$ws.Range("D35").Value = "'11.95"
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("D36").Value = "'0.02298"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "'0.2100"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'0.06044"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "'0.6366"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'4.971"
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "'1.201"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").Value = "'1.386"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'7.804"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'13.28"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").Value = "'0.5906"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'123.15"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").Value = "'1.956"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "'1.147"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").Value = "'0.06852"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'72.19"
$ws.Range("E51").Value = "  -2.28%  "
